# TC08_CDS_Filter_Study-Detection_CCS.xlsx — "startup" sheet lists the
# Neo4j/Web query tabs. Rename the "CasesTab" entry to "ParticipantsTab"
# and leave the selection parked on it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("A2").Select()
